{"js": "// 1) Remove the stray \"_GoBack\" bookmark that Word leaves behind from the\n//    previous editing session (the diff drops <w:bookmarkStart>/<w:bookmarkEnd>\n//    for it from the top of the body).\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2) Insert \", an app available for Android\" right after \"ODK Collect\" and\n//    before the following \". See O...\" sentence in the \"Document Purpose:\"\n//    paragraph. Search on a short, unambiguous snippet so we land on the\n//    first (and only) occurrence of this exact phrase (the words\n//    \"ODK Collect\" alone appear several times later in the document).\nconst anchor = context.document.body.search(\"Ona using ODK Collect\", { matchCase: true });\nanchor.load(\"items\");\nawait context.sync();\n\nif (anchor.items.length === 0) {\n  throw new Error(\"Could not locate the 'Ona using ODK Collect' anchor text\");\n}\n\nanchor.items[0].insertText(\", an app available for Android\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the stray \"_GoBack\" bookmark left over from the previous editing\n#    session (the diff drops <w:bookmarkStart>/<w:bookmarkEnd> for it from\n#    the top of the body).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Insert \", an app available for Android\" right after \"ODK Collect\" and\n#    before the following \". See O...\" sentence in the \"Document Purpose:\"\n#    paragraph. Search on a short, unambiguous snippet so we land on the\n#    first (and only) occurrence of this exact phrase (the words\n#    \"ODK Collect\" alone appear several times later in the document).\n$rng = $d.Content\n$rng.Find.Text = \"Ona using ODK Collect\"\n$found = $rng.Find.Execute()\nif ($found) {\n    $rng.Collapse(0)\n    $rng.InsertAfter(\", an app available for Android\")\n}\n"}
